# Change AddressBook to HealthBook for command sequence diagrams
# (docs/diagrams/DeleteMedicalHistorySequenceDiagram.pptx)
#
# 1. On slide 1, two shapes reference "AddressBook":
#      ":AddressBookParser"       -> ":HealthBookParser"   (only the ":Address" run changes)
#      ":VersionedAddressBook"    -> ":VersionedHealthBook" (only the "VersionedAddressBook" run changes)
# 2. Every "Date Placeholder" (auto date field) on the slide master, the
#    11 slide layouts and the notes master is refreshed from "11/11/18"
#    to "11/12/2018".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide-level text fixes ("Address" -> "Health")
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)

    $hasText = $false
    try {
        $full = $shp.TextFrame.TextRange.Text
        $hasText = $true
    } catch {
        $hasText = $false
    }
    if (-not $hasText) { continue }
    if ($full -notlike "*Address*") { continue }

    $tr2 = $shp.TextFrame2.TextRange

    if ($full -like "*VersionedAddressBook*") {
        # Single run "VersionedAddressBook" (preceded by a separate ":" run).
        # Replace the whole run (character span 2..21) so the run is
        # rewritten in one piece rather than split into several runs.
        $start = $full.IndexOf("VersionedAddressBook") + 1   # 1-based
        $len = "VersionedAddressBook".Length
        $chars = $tr2.Characters($start, $len)
        $chars.Text = "VersionedHealthBook"
    } elseif ($full -like "*Address*") {
        # First paragraph is exactly the run ":Address" (second paragraph
        # holds "BookParser" untouched). Replace that whole run in one go.
        $start = $full.IndexOf(":Address") + 1   # 1-based
        $len = ":Address".Length
        $chars = $tr2.Characters($start, $len)
        $chars.Text = ":Health"
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" placeholder text everywhere
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "11/12/2018"
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout hanging off the slide master
$master = $p.SlideMaster
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes
